# Add two new "Carbon footprint" model parameters to the Setting_Al_cars sheet,
# right after the existing "Segments" (Srsc) parameter row (old row 56 / new row 58
# header "Model flow control" gets pushed down by 2 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting_Al_cars")

# Insert two blank rows where the new parameter rows belong (old row 55 blank
# separator stays put at 55; the two new data rows land at 55 and 56, pushing
# the old gap row + "Model flow control" section, which used to start at row 56,
# down to start at row 58).
$ws.Rows("55:56").Insert()

# Fill in the Descriptor/Name/Index-structure columns in the same order the
# strings were originally typed, so the shared-string table layout matches.
$ws.Range("D55").Value = "Carbon footprint of primary aluminium production"
$ws.Range("C55").Value = "Carbon_Footprint_Primary"
$ws.Range("C56").Value = "Carbon_Footprint_Secondary"
$ws.Range("D56").Value = "Carbon footprint of secondary aluminium production"
$ws.Range("F55").Value = "tS"

# Remaining columns (Version/Aspect order/Layer selection)
$ws.Range("E55").Value = "Carbon_Footprint_Primary"
$ws.Range("G55").Value = "[0,1]"
$ws.Range("H55").Value = "[0]"
$ws.Range("E56").Value = "Carbon_Footprint_Secondary"
$ws.Range("F56").Value = "tS"
$ws.Range("G56").Value = "[0,1]"
$ws.Range("H56").Value = "[0]"

# Match the source workbook's cell formatting for the Parameter_Name (C) and
# Descriptor (D) columns: border on the right of C, themed font in D.
$ws.Range("C54").Copy()
$ws.Range("C55:C56").PasteSpecial(-4122)
$ws.Range("D48").Copy()
$ws.Range("D55:D56").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# The row-insert carries the "Aspect order match" column's border formatting
# down from row 54; the new rows don't use it, so drop it again.
$ws.Range("G55:G56").ClearFormats()
